$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = -0.2896655770681418
$ws.Range("H2").Value = -0.2896655770681418
$ws.Range("I2").Value = -0.3238622076942418
$ws.Range("J2").Value = -0.3238622076942418
$ws.Range("K2").Value = -1.705
$ws.Range("L2").Value = -0.428715111893387
$ws.Range("U2").Value = 2.78
$ws.Range("V2").Value = 0.2291838417147568
$ws.Range("W2").Value = -1.165041765287849
$ws.Range("X2").Value = 0.0779842447179491
$ws.Range("Y2").Value = -1.243026010005798
$ws.Range("Z2").Value = 2.672715053763441
$ws.Range("AA2").Value = 3.001183140718023
$ws.Range("AB2").Value = 0.07685084961221987
$ws.Range("AC2").Value = 2.924332291105803
$ws.Range("AD2").Value = 0.212
$ws.Range("AF2").Value = 0.212
$ws.Range("AG2").Value = -2.568
$ws.Range("AH2").Value = 0.01717711878139686
$ws.Range("AI2").Value = 0.05502206073189722
$ws.Range("AJ2").Value = -0.2685630621208952
$ws.Range("AK2").Value = -2.393289841565702
$ws.Range("AL2").Value = 0.061
$ws.Range("AM2").Value = 0.06
$ws.Range("AN2").Value = -0.2009478672985782
$ws.Range("AO2").Value = -21.11475409836065
$ws.Range("AP2").Value = 2.434123222748815
$ws.Range("AQ2").Value = -21.46666666666667
$ws.Range("B3").Value = "Industrial Solar Holding Europe AB (NGM:ISHE)"
$ws.Range("G3").Value = -8.025641025641026
$ws.Range("H3").Value = -8.025641025641026
$ws.Range("I3").Value = -9.829059829059828
$ws.Range("J3").Value = -9.829059829059828
$ws.Range("K3").Value = -1.2
$ws.Range("L3").Value = -10.25641025641026
$ws.Range("U3").Value = 2.78
$ws.Range("V3").Value = 0.3861111111111111
$ws.Range("W3").Value = -1.938610662358643
$ws.Range("X3").Value = 0.07712332527260096
$ws.Range("Y3").Value = -2.015733987631244
$ws.Range("Z3").Value = -0.6190476190476188
$ws.Range("AA3").Value = 6.084656084656082
$ws.Range("AB3").Value = 0.07693059334984469
$ws.Range("AC3").Value = 6.007725491306237
$ws.Range("AD3").Value = 0.031
$ws.Range("AF3").Value = 0.031
$ws.Range("AG3").Value = -2.749
$ws.Range("AH3").Value = 0.004287097220301479
$ws.Range("AI3").Value = 0.0112278160086925
$ws.Range("AJ3").Value = -0.6176140193215006
$ws.Range("AK3").Value = 144.6842105263182
$ws.Range("AL3").Value = 0.049
$ws.Range("AM3").Value = 0.049
$ws.Range("AN3").Value = -0.0287037037037037
$ws.Range("AO3").Value = -23.46938775510204
$ws.Range("AP3").Value = 2.54537037037037
$ws.Range("AQ3").Value = -23.46938775510204
$ws.Range("B4").Value = "Vadsbo SwitchTech Group AB (publ) (NGM:VADS)"
$ws.Range("G4").Value = -0.05518134715025907
$ws.Range("H4").Value = -0.05518134715025907
$ws.Range("I4").Value = -0.03575129533678757
$ws.Range("J4").Value = -0.03575129533678757
$ws.Range("K4").Value = -0.505
$ws.Range("L4").Value = -0.1308290155440414
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = -0.3914728682170542
$ws.Range("X4").Value = 0.07884516416329723
$ws.Range("Y4").Value = -0.4703180323803515
$ws.Range("Z4").Value = 2.301729278473464
$ws.Range("AA4").Value = -0.08228980322003578
$ws.Range("AB4").Value = 0.07677110587459506
$ws.Range("AC4").Value = -0.1590609090946308
$ws.Range("AD4").Value = 0.181
$ws.Range("AF4").Value = 0.181
$ws.Range("AG4").Value = 0.181
$ws.Range("AH4").Value = 0.03541381334376834
$ws.Range("AI4").Value = 0.1657509157509157
$ws.Range("AJ4").Value = 0.03541381334376834
$ws.Range("AK4").Value = 0.1657509157509157
$ws.Range("AL4").Value = 0.012
$ws.Range("AM4").Value = 0.011
$ws.Range("AN4").Value = 7.239999999999999
$ws.Range("AO4").Value = -11.5
$ws.Range("AP4").Value = 7.239999999999999
$ws.Range("AQ4").Value = -12.54545454545455
